$d = $word.ActiveDocument

# The document contains four "<id>...</id>" markers, each split across three
# runs: "<id>" (Courier New, colored), "p049r_aN" (plain run), "</id>"
# (Courier New, colored). Collapse each trio into a single run reading
# "<id>p049r_N</id>" (dropping the "a"), using the first run's formatting.

$d.Content.Find.Execute("<id>p049r_a1</id>", $false, $false, $false, $false, `
    $false, $true, 1, $false, "<id>p049r_1</id>", 2) | Out-Null

$d.Content.Find.Execute("<id>p049r_a2</id>", $false, $false, $false, $false, `
    $false, $true, 1, $false, "<id>p049r_2</id>", 2) | Out-Null

$d.Content.Find.Execute("<id>p049r_a3</id>", $false, $false, $false, $false, `
    $false, $true, 1, $false, "<id>p049r_3</id>", 2) | Out-Null

$d.Content.Find.Execute("<id>p049r_a4</id>", $false, $false, $false, $false, `
    $false, $true, 1, $false, "<id>p049r_4</id>", 2) | Out-Null
